$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two new rows at the top of the "Zapallo" data block (row 536),
# pushing the existing rows 536:563 down to 538:565.
$ws.Rows("536:537").Insert()

# New row 536: fresh weekly entry (Camote, 1a (guarda))
$ws.Range("A536").Value = 10
$ws.Range("B536").Value = "Vega Modelo de Temuco"
$ws.Range("C536").Value = "La Araucanía"
$ws.Range("D536").Value = 44706
$ws.Range("E536").Value = 9
$ws.Range("F536").Value = 100112045
$ws.Range("G536").Value = "Zapallo"
$ws.Range("H536").Value = "Camote"
$ws.Range("I536").Value = "1a (guarda)"
$ws.Range("J536").Value = 450
$ws.Range("K536").Value = 600
$ws.Range("L536").Value = 700
$ws.Range("M536").Value = 656
$ws.Range("N536").Value = "`$/kilo (volumen en unidades)"
$ws.Range("O536").Value = "Región de O'Higgins"
$ws.Range("P536").Value = 656
$ws.Range("Q536").Value = 1
$ws.Range("R536").Value = "Hortaliza"

# New row 537: fresh weekly entry (Camote, 1a (guarda))
$ws.Range("A537").Value = 10
$ws.Range("B537").Value = "Vega Modelo de Temuco"
$ws.Range("C537").Value = "La Araucanía"
$ws.Range("D537").Value = 44706
$ws.Range("E537").Value = 9
$ws.Range("F537").Value = 100112045
$ws.Range("G537").Value = "Zapallo"
$ws.Range("H537").Value = "Camote"
$ws.Range("I537").Value = "1a (guarda)"
$ws.Range("J537").Value = 700
$ws.Range("K537").Value = 500
$ws.Range("L537").Value = 600
$ws.Range("M537").Value = 543
$ws.Range("N537").Value = "`$/kilo (volumen en unidades)"
$ws.Range("O537").Value = "Región del Maule"
$ws.Range("P537").Value = 543
$ws.Range("Q537").Value = 1
$ws.Range("R537").Value = "Hortaliza"
